$p = $ppt.ActivePresentation

# --- Slide 6: RandomForest "Genauigkeit bei 100 Folds" 0,939 -> 0,93 ---
$s6 = $p.Slides.Item(6)
$tbl6 = $s6.Shapes.Item("Tabelle 5").Table
$tbl6.Cell(2, 2).Shape.TextFrame.TextRange.Text = "0,93"

# --- Slide 7: DecisionTree row values updated ---
$s7 = $p.Slides.Item(7)
$tbl7 = $s7.Shapes.Item("Tabelle 5").Table
$tbl7.Cell(3, 3).Shape.TextFrame.TextRange.Text = "0,79"   # Genauigkeit 0,48 -> 0,79
$tbl7.Cell(3, 4).Shape.TextFrame.TextRange.Text = "0,77"   # Klasse 1   0,33 -> 0,77
$tbl7.Cell(3, 6).Shape.TextFrame.TextRange.Text = "0,68"   # Klasse 3   0,43 -> 0,68
$tbl7.Cell(3, 7).Shape.TextFrame.TextRange.Text = "0,99"   # Klasse 4   1    -> 0,99

# --- Slide 7: caption textbox wording tweak ---
$caption = $s7.Shapes.Item("Textfeld 5")
$caption.TextFrame.TextRange.Text = "-> Klasse 4 wird richtig erkannt"
